$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update threshold values on the remaining rows ---
$ws.Range("B2").Value = 5        # alpha_distance_range Min: 4.5 -> 5
$ws.Range("B3").Value = 4.5      # beta_distance_range Min: 4.6 -> 4.5
$ws.Range("C4").Value = 1.5      # ratio_threshold_range Max: 1.3 -> 1.5
$ws.Range("C6").Value = 15       # pie_threshold_range Max (currently row 6): 20 -> 15

# --- Remove the theta_threshold_range row (row 5) entirely; pie row shifts up to row 5 ---
$ws.Rows(5).Delete()

# --- Reflect the new selection left after editing ---
[void]$ws.Range("B3").Select()

# --- Page setup as recorded by the last save ---
$ps = $ws.PageSetup
$ps.PaperSize = 9      # xlPaperA4
$ps.Orientation = 1    # xlPortrait
